$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows span 2..99. Row 36 is left untouched (data anomaly in source).
# Rows 50..57 get a special update: E -> 10, F -> +10 days (20260201 -> 20260211).
# All other data rows (2..99 excluding 36 and 50..57) get E decremented by 1.

for ($r = 2; $r -le 99; $r++) {
    if ($r -eq 36) {
        continue
    }

    if ($r -ge 50 -and $r -le 57) {
        $ws.Cells.Item($r, 5).Value = 10
        $ws.Cells.Item($r, 6).Value = 20260211
    } else {
        $current = $ws.Cells.Item($r, 5).Value2
        $ws.Cells.Item($r, 5).Value = $current - 1
    }
}
